$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$r = $ws.Range("AC49")
$r.Interior.Pattern = 1
$r.Interior.Color = 0xBFBFBF
Write-Host "Color: " $r.Interior.Color
